$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
$ws.Columns.Item(3).ColumnWidth = 9.94
